$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns E1:P1
$headers = @(
    "contenido_colaborar",
    "contenido_contratación",
    "contenido_investigación",
    "contenido_persoal",
    "contenido_proxecto",
    "contenido_proxectos",
    "titulo_colaborar",
    "titulo_contratación",
    "titulo_investigación",
    "titulo_persoal",
    "titulo_proxecto",
    "titulo_proxectos"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 5 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Row 4: UVigoProfesor (re-scrape, newer timestamp), E4:P4 left blank
$ws.Range("A4").Value = "UVigoProfesor"
$ws.Range("B4").Value = "NO"
$ws.Range("C4").Value = "2025-10-09 16:29:09"
$ws.Range("D4").Value = "https://secretaria.uvigo.gal/uv/web/convocatoria/public/index"

# Row 5: USCEmprego (re-scrape, newer timestamp), E5:P5 left blank
$ws.Range("A5").Value = "USCEmprego"
$ws.Range("B5").Value = "ERROR"
$ws.Range("C5").Value = "2025-10-09 16:29:09"
$ws.Range("D5").Value = "https://www.usc.gal/gl/emprego"

# Row 6: new USCConvocatorias entry with content/title counts
$ws.Range("A6").Value = "USCConvocatorias"
$ws.Range("C6").Value = "2025-10-09 16:29:09"
$ws.Range("D6").Value = "https://www.usc.gal/gl/investigar-na-usc/convocatorias"

$counts = @(12, 10, 12, 12, 12, 12, 12, 10, 12, 12, 12, 12)
for ($i = 0; $i -lt $counts.Length; $i++) {
    $col = 5 + $i
    $ws.Cells.Item(6, $col).Value = $counts[$i]
}
